$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 62.25498566666667
$ws.Range("H2").Value = 186.764957
$ws.Range("I2").Value = 0.7607543248383142
$ws.Range("J2").Value = 0.7607543248383141
$ws.Range("M2").Value = 62.25498566666667
$ws.Range("N2").Value = 186.764957
$ws.Range("O2").Value = 0.7607543248383142
$ws.Range("P2").Value = 0.7607543248383141
$ws.Range("Q2").Value = 3875.683240356872
$ws.Range("R2").Value = 34881.14916321186
$ws.Range("S2").Value = 0.5787471427601992
$ws.Range("T2").Value = 0.5787471427601991
$ws.Range("G3").Value = 62.25498566666667
$ws.Range("H3").Value = 186.764957
$ws.Range("I3").Value = 0.7607543248383142
$ws.Range("J3").Value = 0.7607543248383141
$ws.Range("O3").Value = 0.1681603168407971
$ws.Range("P3").Value = 0.1681603168407971
$ws.Range("Q3").Value = 856.6972285191999
$ws.Range("R3").Value = 7710.275056672799
$ws.Range("S3").Value = 0.1279286883028176
$ws.Range("T3").Value = 0.1279286883028176
$ws.Range("G4").Value = 62.25498566666667
$ws.Range("H4").Value = 186.764957
$ws.Range("I4").Value = 0.7607543248383142
$ws.Range("J4").Value = 0.7607543248383141
$ws.Range("M4").Value = 5.817144666666667
$ws.Range("N4").Value = 17.451434
$ws.Range("O4").Value = 0.07108535832088886
$ws.Range("P4").Value = 0.07108535832088884
$ws.Range("Q4").Value = 362.1462578442598
$ws.Range("R4").Value = 3259.316320598338
$ws.Range("S4").Value = 0.05407849377529744
$ws.Range("T4").Value = 0.05407849377529742
$ws.Range("I5").Value = 0.1681603168407971
$ws.Range("J5").Value = 0.1681603168407971
$ws.Range("M5").Value = 62.25498566666667
$ws.Range("N5").Value = 186.764957
$ws.Range("O5").Value = 0.7607543248383142
$ws.Range("P5").Value = 0.7607543248383141
$ws.Range("Q5").Value = 856.6972285191999
$ws.Range("R5").Value = 7710.275056672799
$ws.Range("S5").Value = 0.1279286883028176
$ws.Range("T5").Value = 0.1279286883028176
$ws.Range("I6").Value = 0.1681603168407971
$ws.Range("J6").Value = 0.1681603168407971
$ws.Range("O6").Value = 0.1681603168407971
$ws.Range("P6").Value = 0.1681603168407971
$ws.Range("S6").Value = 0.02827789215999727
$ws.Range("T6").Value = 0.02827789215999726
$ws.Range("I7").Value = 0.1681603168407971
$ws.Range("J7").Value = 0.1681603168407971
$ws.Range("M7").Value = 5.817144666666667
$ws.Range("N7").Value = 17.451434
$ws.Range("O7").Value = 0.07108535832088886
$ws.Range("P7").Value = 0.07108535832088884
$ws.Range("Q7").Value = 80.05032304580422
$ws.Range("R7").Value = 720.452907412238
$ws.Range("S7").Value = 0.01195373637798226
$ws.Range("T7").Value = 0.01195373637798226
$ws.Range("G8").Value = 5.817144666666667
$ws.Range("H8").Value = 17.451434
$ws.Range("I8").Value = 0.07108535832088886
$ws.Range("J8").Value = 0.07108535832088884
$ws.Range("M8").Value = 62.25498566666667
$ws.Range("N8").Value = 186.764957
$ws.Range("O8").Value = 0.7607543248383142
$ws.Range("P8").Value = 0.7607543248383141
$ws.Range("Q8").Value = 362.1462578442598
$ws.Range("R8").Value = 3259.316320598338
$ws.Range("S8").Value = 0.05407849377529744
$ws.Range("T8").Value = 0.05407849377529742
$ws.Range("G9").Value = 5.817144666666667
$ws.Range("H9").Value = 17.451434
$ws.Range("I9").Value = 0.07108535832088886
$ws.Range("J9").Value = 0.07108535832088884
$ws.Range("O9").Value = 0.1681603168407971
$ws.Range("P9").Value = 0.1681603168407971
$ws.Range("Q9").Value = 80.05032304580422
$ws.Range("R9").Value = 720.452907412238
$ws.Range("S9").Value = 0.01195373637798226
$ws.Range("T9").Value = 0.01195373637798226
$ws.Range("G10").Value = 5.817144666666667
$ws.Range("H10").Value = 17.451434
$ws.Range("I10").Value = 0.07108535832088886
$ws.Range("J10").Value = 0.07108535832088884
$ws.Range("M10").Value = 5.817144666666667
$ws.Range("N10").Value = 17.451434
$ws.Range("O10").Value = 0.07108535832088886
$ws.Range("P10").Value = 0.07108535832088884
$ws.Range("Q10").Value = 33.83917207292844
$ws.Range("R10").Value = 304.552548656356
$ws.Range("S10").Value = 0.005053128167609163
$ws.Range("T10").Value = 0.005053128167609161
